# Update the threshold values in the 2His_1Asp sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (alpha_distance_range): Min 4 -> 3.8, Max 12 -> 12.2
$ws.Range("B2").Value = 3.8
$ws.Range("C2").Value = 12.2

# Row 3 (beta_distance_range): Max 10.5 -> 11.1
$ws.Range("C3").Value = 11.1

# Row 4 (ratio_threshold_range): Min 0.75 -> 0.7, Max 1.5 -> 1.65
$ws.Range("B4").Value = 0.7
$ws.Range("C4").Value = 1.65

# Row 5 (pie_threshold_range): Max 25 -> 28
$ws.Range("C5").Value = 28
